# Weekly update: a new price record for the week of 2023-12-20 is added
# for "Vega Monumental Concepción - Poroto verde". The new record is
# inserted as row 36 (rows are kept in reverse-chronological / most
# recent-entry-first order for this subset), pushing the previously
# existing rows 36-97 down to rows 37-98.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 36; this shifts rows 36..97 down to 37..98
# and keeps all existing formatting (e.g. the date style on column D).
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row 36 with the new weekly record.
$ws.Range("A36").Value = 11
$ws.Range("B36").Value = "Vega Monumental Concepción"
$ws.Range("C36").Value = "Bíobío"
$ws.Range("D36").Value = 45280
$ws.Range("E36").Value = 8
$ws.Range("F36").Value = 100112031
$ws.Range("G36").Value = "Poroto verde"
$ws.Range("H36").Value = "Magnum"
$ws.Range("I36").Value = "Primera"
$ws.Range("J36").Value = 100
$ws.Range("K36").Value = 21000
$ws.Range("L36").Value = 22000
$ws.Range("M36").Value = 21500
$ws.Range("N36").Value = "$/saco 25 kilos"
$ws.Range("O36").Value = "Región Metropolitana"
$ws.Range("P36").Value = 860
$ws.Range("Q36").Value = 25
$ws.Range("R36").Value = "Hortaliza"
